$d = $word.ActiveDocument

# Step 1: Insert a new empty paragraph right after the "Introduction" heading
# paragraph, in the correct structural position (before the chapter-1
# bookmarkEnd), by using Find/Replace to append a paragraph break + a
# placeholder marker right after the word "Introduction".
$rng = $d.Content
$found = $rng.Find.Execute("Introduction", $false, $false, $false, $false, $false, $true, 1, $false, "Introduction^pQUOTE_PLACEHOLDER", 2)
Write-Host "Found1: $found"

# Step 2: Locate the newly created placeholder paragraph (paragraph 6) and
# replace its contents with the quote, built from explicit OOXML so that
# each differently-formatted chunk of text stays in its own <w:r> run
# (mirroring how the text would have been typed originally), while
# preserving the BlockText paragraph style.
$quotePara = $d.Paragraphs.Item(6)
$quoteRange = $quotePara.Range
Write-Host "quoteRange text=[$($quoteRange.Text)] Start=$($quoteRange.Start) End=$($quoteRange.End)"

$xmlFrag = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:pStyle w:val="BlockText"/></w:pPr>
<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">&#x201C;</w:t></w:r>
<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">My data is everywhere, and I am nowhere.</w:t></w:r>
<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">&#x201D;</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t xml:space="preserve">&#x2013; Imogen Heap, musician and digital rights advocate, speaking at MyData 2019.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$quoteRange.InsertXML($xmlFrag)

Write-Host "Paragraph count: $($d.Paragraphs.Count)"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Host "$i`: [$($p.Range.Text)]"
}
